$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').Value = '21.154.45'
$ws.Range('E2').Value = '  +3.58%  '

# Row 3
$ws.Range('D3').Value = '1.540.70'
$ws.Range('E3').Value = '  +5.10%  '

# Row 4
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.011'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.42%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.9671'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.01%  '

# Row 6
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '282.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +2.75%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.3633'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  -0.39%  '

# Row 8
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.3197'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.31%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '40.83'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.73%  '

# Row 10
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.100'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.40%  '

# Row 11
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.06829'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +3.46%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.004'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.51%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.695'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +5.25%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '18.76'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +3.40%  '

# Row 15
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '6.376'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +3.85%  '

# Row 16
$ws.Range('B16').Value = 'ShibaInu'
$ws.Range('C16').Value = 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib'
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '0.00001045'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +2.11%  '

# Row 17
$ws.Range('B17').Value = 'Dai'
$ws.Range('C17').Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.9676'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  -0.90%  '

# Row 18
$ws.Range('D18').Value = '1.537.39'
$ws.Range('E18').Value = '  +4.92%  '

# Row 19
$ws.Range('E19').Value = '  +3.86%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '72.68'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +4.54%  '

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '5.722'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +5.37%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '15.09'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +5.05%  '

# Row 23
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '11.37'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +4.08%  '

# Row 24
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '2.335'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +4.23%  '

# Row 25
$ws.Range('D25').Value = '21.192.36'
$ws.Range('E25').Value = '  +3.65%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '148.76'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.32%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.227'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +6.28%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '17.69'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.03%  '

# Row 29
$ws.Range('D29').Value = '1.708.18'
$ws.Range('E29').Value = '  +5.50%  '

# Row 30
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '118.58'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +4.23%  '

# Row 31
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.012'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.57%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.8583'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +7.52%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.230'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +5.42%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.08019'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.99%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.504'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -0.09%  '

# Row 36
$ws.Range('E36').Value = '  +5.42%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '4.937'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +5.08%  '

# Row 38
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.05863'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +2.42%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.02109'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +3.55%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '10.68'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +3.21%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '7.761'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.51%  '

# Row 42
$ws.Range('B42').Value = 'Algorand'
$ws.Range('C42').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.1921'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +3.24%  '

# Row 43
$ws.Range('B43').Value = 'Frax'
$ws.Range('C43').Value = 'https://coinranking.com/coin/KfWtaeV1W+frax-frax'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.9679'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.25%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.5460'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.54%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '12.59'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +5.30%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.570'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.55%  '

# Row 47
$ws.Range('B47').Value = 'Decentraland'
$ws.Range('C47').Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.5447'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +5.27%  '

# Row 48
$ws.Range('B48').Value = 'Quant'
$ws.Range('C48').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '122.15'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +4.30%  '

# Row 49
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.875'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +6.55%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.06586'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +2.62%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.9914'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +0.18%  '
